# Swap the two embedded themes: the "Integral" / Red-Violet palette
# (currently driving the slide master via theme2.xml) and the generic
# "Office Theme" palette (currently only used by the notes master via
# theme1.xml) trade places, per:
#   theme1.xml: Office Theme  -> Integral / Red Violet
#   theme2.xml: Integral      -> Office Theme
#
# The exposed object model only lets us reach the live theme through
# SlideMaster / NotesMaster / HandoutMaster.Theme - they all resolve to
# the same underlying theme object, which corresponds to theme2.xml on
# save. We drive that object's 12 scheme colors to the "Office Theme"
# palette so theme2.xml ends up matching the target content.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$cs = $master.Theme.ThemeColorScheme

# Order is dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (COM RGB = 0xBBGGRR)
$cs.Item(1).RGB  = 0        # dk1      000000
$cs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388  # dk2      44546A
$cs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501  # accent2  ED7D31
$cs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$cs.Item(8).RGB  = 49407    # accent4  FFC000
$cs.Item(9).RGB  = 12874308 # accent5  4472C4
$cs.Item(10).RGB = 4697456  # accent6  70AD47
$cs.Item(11).RGB = 12673797 # hlink    0563C1
$cs.Item(12).RGB = 7491477  # folHlink 954F72
